$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # 展览
$ws1.Range("F3").Value = 3322
$ws1.Range("F7").Value = 3827
$ws1.Range("F11").Value = 8559
$ws1.Range("F12").Value = 8559
$ws1.Range("F16").Value = 104
$ws1.Range("F17").Value = 326
$ws1.Range("F19").Value = 82
$ws1.Range("F21").Value = 352
$ws1.Range("F22").Value = 10811
$ws1.Range("F23").Value = 10811
$ws1.Range("F32").Value = 2671
$ws1.Range("F35").Value = 2075
$ws1.Range("F37").Value = 41
$ws1.Range("F40").Value = 4075
$ws1.Range("F41").Value = 2230
$ws1.Range("F44").Value = 3024
$ws1.Range("F45").Value = 1239
$ws1.Range("F47").Value = 751
$ws1.Range("F48").Value = 339
$ws1.Range("F49").Value = 310
$ws1.Range("F50").Value = 41
$ws1.Range("F51").Value = 122

$ws2 = $wb.Worksheets.Item(2)  # 演出
$ws2.Range("F9").Value = 6
$ws2.Range("F16").Value = 10

$ws3 = $wb.Worksheets.Item(3)  # 本地生活
$ws3.Range("F3").Value = 24

$ws4 = $wb.Worksheets.Item(4)  # 全部类型
$ws4.Range("F5").Value = 3322
$ws4.Range("F10").Value = 3827
$ws4.Range("F16").Value = 8559
$ws4.Range("F20").Value = 104
$ws4.Range("F21").Value = 326
$ws4.Range("F23").Value = 82
$ws4.Range("F25").Value = 10811
$ws4.Range("F28").Value = 24
$ws4.Range("F35").Value = 2671
$ws4.Range("F38").Value = 2075
$ws4.Range("F40").Value = 41
$ws4.Range("F44").Value = 2230
$ws4.Range("F45").Value = 3024
$ws4.Range("F47").Value = 1239
$ws4.Range("F48").Value = 339
$ws4.Range("F49").Value = 310
$ws4.Range("F50").Value = 41
$ws4.Range("F51").Value = 122
